# Applies the cryptos.xlsx price/volume/coin-swap refresh described in the
# commit "Updated cryptos list on Wed Nov 29 12:47:44 UTC 2023 with GitHub
# Actions" - per-row Price (D) and Volume(1h) (E) updates, plus two coin-pair
# re-sorts (Polkadot/Polygon rows 15-16, Monero/Cosmos rows 26-27,
# InjectiveProtocol/Aave rows 41-42) where Name/Link/Price/Volume all swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.267.04'
$ws.Range("E2").Value = '  +2.67%  '
$ws.Range("D3").Value = '2.058.00'
$ws.Range("E3").Value = '  +1.60%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''229.18'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("D7").Value = '''61.12'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +8.76%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '''0.386'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("D10").Value = '''0.0830'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +5.95%  '
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").Value = '''14.83'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.70%  '
$ws.Range("D13").Value = '2.365.01'
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("D14").Value = '''21.19'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.71%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.762'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '''5.34'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.70%  '
$ws.Range("D17").Value = '2.058.92'
$ws.Range("E17").Value = '  +1.39%  '
$ws.Range("D18").Value = '38.217.78'
$ws.Range("E18").Value = '  +2.69%  '
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("D20").Value = '''69.83'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("D21").Value = '0.0₃0837'
$ws.Range("E21").Value = '  +2.23%  '
$ws.Range("D22").Value = '''225.51'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '''2.43'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("D25").Value = '''2.23'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''9.27'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.37%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''166.12'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.49%  '
$ws.Range("E28").Value = '  +2.50%  '
$ws.Range("D29").Value = '''18.97'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.34%  '
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").Value = '''0.121'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.54%  '
$ws.Range("D32").Value = '''4.52'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.36%  '
$ws.Range("D33").Value = '''4.60'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.34%  '
$ws.Range("E34").Value = '  +2.29%  '
$ws.Range("D35").Value = '''0.0606'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").Value = '''6.38'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +14.93%  '
$ws.Range("D37").Value = '''2.29'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.60%  '
$ws.Range("D38").Value = '''3.28'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.18%  '
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").Value = '1.527.67'
$ws.Range("E40").Value = '  +3.99%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '''97.88'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.73%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = '''16.99'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.55%  '
$ws.Range("D43").Value = '''0.0216'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.30%  '
$ws.Range("E44").Value = '  +1.65%  '
$ws.Range("D45").Value = '''0.0932'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.12%  '
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("D47").Value = '''4.02'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -10.22%  '
$ws.Range("D48").Value = '''1.02'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.62%  '
$ws.Range("D49").Value = '''2.99'
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = '''7.04'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.57%  '
$ws.Range("D51").Value = '2.252.72'
$ws.Range("E51").Value = '  +1.78%  '
